$d = $word.ActiveDocument

$d.Content.Find.Execute("45×54=2430", $true, $false, $false, $false, $false, $true, 1, $false, "13×81=1053", 2) | Out-Null
$d.Content.Find.Execute("64×18=1152", $true, $false, $false, $false, $false, $true, 1, $false, "92×38=3496", 2) | Out-Null
$d.Content.Find.Execute("57×98=5586", $true, $false, $false, $false, $false, $true, 1, $false, "62×81=5022", 2) | Out-Null
$d.Content.Find.Execute("33×44=1452", $true, $false, $false, $false, $false, $true, 1, $false, "92×83=7636", 2) | Out-Null
$d.Content.Find.Execute("59×57=3363", $true, $false, $false, $false, $false, $true, 1, $false, "42×56=2352", 2) | Out-Null
$d.Content.Find.Execute("40×24=960", $true, $false, $false, $false, $false, $true, 1, $false, "85×56=4760", 2) | Out-Null
$d.Content.Find.Execute("50×37=1850", $true, $false, $false, $false, $false, $true, 1, $false, "44×18=792", 2) | Out-Null
$d.Content.Find.Execute("88×53=4664", $true, $false, $false, $false, $false, $true, 1, $false, "43×71=3053", 2) | Out-Null
$d.Content.Find.Execute("28×92=2576", $true, $false, $false, $false, $false, $true, 1, $false, "32×94=3008", 2) | Out-Null
$d.Content.Find.Execute("60×82=4920", $true, $false, $false, $false, $false, $true, 1, $false, "41×63=2583", 2) | Out-Null
$d.Content.Find.Execute("12×71=852", $true, $false, $false, $false, $false, $true, 1, $false, "17×85=1445", 2) | Out-Null
$d.Content.Find.Execute("95×54=5130", $true, $false, $false, $false, $false, $true, 1, $false, "51×87=4437", 2) | Out-Null
$d.Content.Find.Execute("25×69=1725", $true, $false, $false, $false, $false, $true, 1, $false, "86×14=1204", 2) | Out-Null
$d.Content.Find.Execute("75×74=5550", $true, $false, $false, $false, $false, $true, 1, $false, "99×78=7722", 2) | Out-Null
$d.Content.Find.Execute("19×18=342", $true, $false, $false, $false, $false, $true, 1, $false, "29×56=1624", 2) | Out-Null
$d.Content.Find.Execute("73×27=1971", $true, $false, $false, $false, $false, $true, 1, $false, "50×86=4300", 2) | Out-Null
$d.Content.Find.Execute("64×19=1216", $true, $false, $false, $false, $false, $true, 1, $false, "44×61=2684", 2) | Out-Null
$d.Content.Find.Execute("80×49=3920", $true, $false, $false, $false, $false, $true, 1, $false, "35×44=1540", 2) | Out-Null
$d.Content.Find.Execute("25×77=1925", $true, $false, $false, $false, $false, $true, 1, $false, "65×83=5395", 2) | Out-Null
$d.Content.Find.Execute("63×61=3843", $true, $false, $false, $false, $false, $true, 1, $false, "19×86=1634", 2) | Out-Null
$d.Content.Find.Execute("74×84=6216", $true, $false, $false, $false, $false, $true, 1, $false, "19×84=1596", 2) | Out-Null
$d.Content.Find.Execute("68×69=4692", $true, $false, $false, $false, $false, $true, 1, $false, "56×63=3528", 2) | Out-Null
$d.Content.Find.Execute("24×85=2040", $true, $false, $false, $false, $false, $true, 1, $false, "26×32=832", 2) | Out-Null
$d.Content.Find.Execute("74×62=4588", $true, $false, $false, $false, $false, $true, 1, $false, "43×85=3655", 2) | Out-Null
$d.Content.Find.Execute("11×16=176", $true, $false, $false, $false, $false, $true, 1, $false, "17×48=816", 2) | Out-Null
